$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 61, shifting existing rows 61.. down by 2.
$ws.Range("A61:A62").EntireRow.Insert()

# New row 61: Membrillo / Primera, D=45036
$ws.Range("A61").Value = 6
$ws.Range("B61").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C61").Value = "Metropolitana"
$ws.Range("D61").Value = 45036
$ws.Range("E61").Value = 13
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100104
$ws.Range("H61").Value = "Frutos de pepita"
$ws.Range("I61").Value = 100104003
$ws.Range("J61").Value = "Membrillo"
$ws.Range("K61").Value = "Champion"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 6
$ws.Range("N61").Value = 230000
$ws.Range("O61").Value = 230000
$ws.Range("P61").Value = 230000
$ws.Range("Q61").Value = "$/bins (450 kilos)"
$ws.Range("R61").Value = "Región de O'Higgins"
$ws.Range("S61").Value = 511
$ws.Range("T61").Value = 450

# New row 62: Membrillo / Segunda, D=45036
$ws.Range("A62").Value = 6
$ws.Range("B62").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C62").Value = "Metropolitana"
$ws.Range("D62").Value = 45036
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100104
$ws.Range("H62").Value = "Frutos de pepita"
$ws.Range("I62").Value = 100104003
$ws.Range("J62").Value = "Membrillo"
$ws.Range("K62").Value = "Champion"
$ws.Range("L62").Value = "Segunda"
$ws.Range("M62").Value = 10
$ws.Range("N62").Value = 200000
$ws.Range("O62").Value = 200000
$ws.Range("P62").Value = 200000
$ws.Range("Q62").Value = "$/bins (450 kilos)"
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 444
$ws.Range("T62").Value = 450
